# Append two new data rows (57 and 58) to Sheet1, mirroring the existing
# two-row pattern already present in rows 53/54 (and repeated at 55/56).
#
# Row 57 : datetime / "Andy McAllister" / date / "HOCKERWOOD F/R - F614 - LF1949"
# Row 58 : datetime / "sfdlkhsdkj" / "dkjhgfskjfhk" / "dskjfhskdjhfksj"
#
# The date/datetime number formats must be set using the exact escaped
# format-code strings already used by the workbook's existing styles
# (yyyy\-mm\-dd\ hh:mm:ss / dd\-mm\-yyyy) so the engine reuses the current
# style indices instead of minting duplicate numFmt/cellXf entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 57
$ws.Cells.Item(57, 1).Value = 43424.571840277778
$ws.Cells.Item(57, 1).NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Cells.Item(57, 2).Value = "Andy McAllister"
$ws.Cells.Item(57, 3).Value = 43424
$ws.Cells.Item(57, 3).NumberFormat = "dd\-mm\-yyyy"
$ws.Cells.Item(57, 4).Value = "HOCKERWOOD F/R - F614 - LF1949"

# Row 58
$ws.Cells.Item(58, 1).Value = 43424.571840277778
$ws.Cells.Item(58, 1).NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Cells.Item(58, 2).Value = "sfdlkhsdkj"
$ws.Cells.Item(58, 3).Value = "dkjhgfskjfhk"
$ws.Cells.Item(58, 4).Value = "dskjfhskdjhfksj"
